$d = $word.ActiveDocument

# Helper: find a paragraph containing the given unique text and delete the
# whole paragraph (including its end-of-paragraph mark).
function Remove-ParagraphContaining($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1)
        $para.Range.Delete()
    }
}

# Drop the two "Do poprawienia" bullet items that were resolved/removed.
Remove-ParagraphContaining("Destination – Create")

# Drop the two "Do dodania" bullet items that were moved to the backend /
# no longer apply.
Remove-ParagraphContaining("Resource Type – List")
Remove-ParagraphContaining("Distance – dodanie nowej tabeli")
Remove-ParagraphContaining("Route – distance – dodanie nowego atrybutu")

# Split the "Zrobienie testów" run in two (right after "Jasmine oraz ") and
# re-insert the (previously removed) _GoBack bookmark at that split point.
$marker = $d.Content
$marker.Find.Execute("Jasmine oraz Karma", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$splitAt = $marker.Start + "Jasmine oraz ".Length
$bmRange = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $bmRange)
